$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header texts (row 1)
$ws.Range("A1").Value = "County"
$ws.Range("B1").Value = "County code"
$ws.Range("C1").Value = "Monitoring Sites"
$ws.Range("D1").Value = "Data Samples"
$ws.Range("E1").Value = "Highest Fecal Coliform Concentration"
$ws.Range("F1").Value = "Units Measured"

# Header formatting
$ws.Range("A1:F1").Font.Bold = $true
$ws.Range("E1").WrapText = $true
$ws.Rows("1").RowHeight = 34

# Column widths (char units -> raw xml width = input + 5/6, rounded to 1/6)
$ws.Columns("B").ColumnWidth = 12.1666667
$ws.Columns("C:D").ColumnWidth = 14.6666667
$ws.Columns("E").ColumnWidth = 21
$ws.Columns("F").ColumnWidth = 18.1666667

# Data cell alignment (center)
$ws.Range("A2:F7").HorizontalAlignment = -4108

# Selection
$excel.Goto($ws.Range("A2:F7"))
